# v3.0 Closed reviewer verification for Login , Delete user and CATEGORIES features
#
# This script applies the substantive content edits described by the commit:
#   1. On REVIEW-SHEET, mark the "Reviewer verification" (column J) status as
#      "closed" for the three LOGIN review rows (30-32) and the DELETEUSER
#      review row (33) -- these were "open" before.
#   2. On VERSION-HISTORY, append a new v3.0 entry describing the change.

$wb = $excel.ActiveWorkbook

# --- 1. REVIEW-SHEET: close out reviewer verification for Login & Delete user rows ---
$reviewSheet = $wb.Worksheets.Item("REVIEW-SHEET")
$reviewSheet.Range("J30:J33").Value = "closed"
# Touch the font so Excel commits fresh, explicit direct formatting on these
# cells (matching how the list-validated cells looked after being
# re-confirmed as "closed" by the reviewer).
$reviewSheet.Range("J30:J33").Font.ColorIndex = 1

# --- 2. VERSION-HISTORY: append the v3.0 changelog entry ---
$historySheet = $wb.Worksheets.Item("VERSION-HISTORY")

# Copy the formatting of the last existing entry (row 21) down onto the new
# row 22 so the new row visually matches the rest of the table.
$historySheet.Range("A21:D21").Copy()
$historySheet.Range("A22").PasteSpecial(-4122)

$historySheet.Range("A22").Value = "v3.0"
$historySheet.Range("B22").Value = "Hala Eldaly"
$historySheet.Range("C22").Value = "Closed reviewer verification for Login , Delete user and CATEGORIES features"
$historySheet.Range("D22").Value = 45787
